# edit.ps1 - applies the "IT Experts" rebranding edit to the questionnaire.
#
# Summary of the change (per the unified diff):
#   1. Title paragraph: "Respondent's Profile Questionnaire"
#        -> " " + "IT Experts" + " " + "Profile Questionnaire"  (4 runs, no proofErr)
#   2. Salutation paragraph: "Dear Respondents,"
#        -> "Dear " + "IT Experts" + ","                         (3 runs, no proofErr)
#   3. Section heading paragraph: "Profile of the respondents:"
#        -> "Profile of the " + "IT Experts" + ":"                (3 runs)
#
# The underlying runtime merges adjacent same-format runs whenever a Find/
# Replace or Range.Text edit crosses a run/proofErr boundary, which makes it
# impossible to reproduce the exact run layout (and proofErr removal) the
# diff wants using plain text replacement. Instead we rebuild the affected
# paragraph's contents explicitly via Range.InsertXML, preserving the
# paragraph's own <w:pPr> (read back dynamically so we do not have to hard
# code rsid/paraId values) and emitting exactly the run list we want.

$d = $word.ActiveDocument

function Get-ParagraphOpenTag {
    # Returns the paragraph's opening "<w:p ...>...</w:pPr>" (or just the
    # bare "<w:p ...>" / "<w:p>" when the paragraph carries no pPr) so the
    # replacement fragment keeps the paragraph's own formatting/identity.
    param($range)
    $xml = $range.WordOpenXML
    $startIdx = $xml.IndexOf("<w:p ")
    if ($startIdx -lt 0) { $startIdx = $xml.IndexOf("<w:p>") }
    $pPrEndIdx = $xml.IndexOf("</w:pPr>", $startIdx)
    if ($pPrEndIdx -lt 0) {
        $tagEnd = $xml.IndexOf(">", $startIdx)
        return $xml.Substring($startIdx, $tagEnd + 1 - $startIdx)
    }
    return $xml.Substring($startIdx, ($pPrEndIdx + 8) - $startIdx)
}

function Set-ParagraphRuns {
    # Replaces a whole paragraph's runs with $runsXml, keeping its own
    # <w:pPr> intact, and without leaving any <w:proofErr/> markers behind.
    param($paragraph, [string]$runsXml)

    $range = $paragraph.Range
    $openTag = Get-ParagraphOpenTag $range

    $fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $openTag + $runsXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $range.InsertXML($fragment)
}

# Shared run-properties blocks (copied verbatim from the surrounding runs so
# the new/edited runs keep identical character formatting).
$rPrTitle  = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:lang w:val="fr-FR"/></w:rPr>'
$rPrHeading = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr>'

# ---------------------------------------------------------------------
# 1) Title paragraph: "<space>Respondent's Profile Questionnaire"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Respondent*Profile Questionnaire*") {
        $runs = ''
        $runs += '<w:r>' + $rPrTitle + '<w:t xml:space="preserve"> </w:t></w:r>'
        $runs += '<w:r>' + $rPrTitle + '<w:t>IT Experts</w:t></w:r>'
        $runs += '<w:r>' + $rPrTitle + '<w:t xml:space="preserve"> </w:t></w:r>'
        $runs += '<w:r>' + $rPrTitle + '<w:t>Profile Questionnaire</w:t></w:r>'
        Set-ParagraphRuns $p $runs
        break
    }
}

# ---------------------------------------------------------------------
# 2) Salutation paragraph: "Dear Respondents,"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Dear Respondents*") {
        $runs = ''
        $runs += '<w:r>' + $rPrTitle + '<w:t xml:space="preserve">Dear </w:t></w:r>'
        $runs += '<w:r>' + $rPrTitle + '<w:t>IT Experts</w:t></w:r>'
        $runs += '<w:r>' + $rPrTitle + '<w:t>,</w:t></w:r>'
        Set-ParagraphRuns $p $runs
        break
    }
}

# ---------------------------------------------------------------------
# 3) Section heading paragraph: "Profile of the respondents:"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Profile of the respondents*") {
        $runs = ''
        $runs += '<w:r>' + $rPrHeading + '<w:t xml:space="preserve">Profile of the </w:t></w:r>'
        $runs += '<w:r>' + $rPrHeading + '<w:t>IT Experts</w:t></w:r>'
        $runs += '<w:r>' + $rPrHeading + '<w:t>:</w:t></w:r>'
        Set-ParagraphRuns $p $runs
        break
    }
}
